# The commit swaps the puzzle's starting board (A1:D4 on the "Manhattan"
# sheet) for a different scramble of the same 4x4 15-puzzle. Everything
# else on both sheets (Manhattan + Linear) is a formula that recalculates
# automatically from these four input cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Manhattan")
$ws.Activate()

$ws.Range("A1").Value = 14
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 8
$ws.Range("D1").Value = 0

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 12
$ws.Range("D2").Value = 10

$ws.Range("A3").Value = 11
$ws.Range("B3").Value = 13
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 2

$ws.Range("A4").Value = 7
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 9

# Move the selection/active cell to match the saved view state.
$ws.Range("E14").Select()
